$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B4").Value = 2
$ws.Range("B58").Value = 3
$ws.Range("B59").Value = 2
$ws.Range("B60").Value = 1
